$wb = $excel.ActiveWorkbook

# Update status text from "Ready for handoff" to "In Translation" on each sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Set the affected columns to reflect the shorter status text's autofit width.
# (Target stored width is 13.4101845877511; the host quantizes ColumnWidth
# assignments to an internal 1/6-wide pixel grid, so 12.5 is the closest input
# that lands on the nearest achievable grid point, 13.333333333333334.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
